$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (prices & volume change %).
$updates = @{
    'D2' = '96.081.43'
    'E2' = '  +0.52%  '
    'D3' = '3.546.79'
    'E3' = '  -1.46%  '
    'E4' = '  -0.08%  '
    'D5' = '239.83'
    'E5' = '  +0.49%  '
    'D6' = '651.10'
    'E6' = '  -0.35%  '
    'D7' = '1.63'
    'E7' = '  +10.49%  '
    'D8' = '0.405'
    'E8' = '  +0.12%  '
    'E9' = '  +6.65%  '
    'E10' = '  +0.02%  '
    'D11' = '3.543.76'
    'E11' = '  -1.59%  '
    'D12' = '43.33'
    'E12' = '  +1.15%  '
    'E13' = '  +0.70%  '
    'D14' = '6.36'
    'E14' = '  +0.53%  '
    'D15' = '4.204.94'
    'E15' = '  -1.58%  '
    'D16' = '96.013.86'
    'E16' = '  +0.54%  '
    'D17' = '0.0000259'
    'E17' = '  +1.61%  '
    'D18' = '3.553.74'
    'E18' = '  -1.50%  '
    'D19' = '7.83'
    'E19' = '  -0.39%  '
    'D20' = '12.37'
    'E20' = '  -2.38%  '
    'E21' = '  -1.96%  '
    'D22' = '0.529'
    'E22' = '  +8.08%  '
    'D23' = '504.78'
    'E23' = '  -1.00%  '
    'E24' = '  -6.32%  '
    'D25' = '6.85'
    'E25' = '  +3.44%  '
    'D26' = '0.0000197'
    'E26' = '  +0.63%  '
    'D27' = '95.73'
    'E27' = '  -1.14%  '
    'D28' = '12.68'
    'E28' = '  -0.31%  '
    'D29' = '3.737.20'
    'E29' = '  -1.68%  '
    'E30' = '  +7.51%  '
    'E31' = '  -4.30%  '
    'E32' = '  -0.31%  '
    'D33' = '0.999'
    'E33' = '  +0.02%  '
    'E34' = '  +2.13%  '
    'E35' = '  +0.49%  '
    'D36' = '31.16'
    'E36' = '  -2.53%  '
    'D37' = '8.69'
    'E37' = '  +5.75%  '
    'D38' = '610.44'
    'E38' = '  +6.27%  '
    'E39' = '  -0.06%  '
    'D40' = '1.60'
    'E40' = '  +7.64%  '
    'E41' = '  +0.05%  '
    'D42' = '0.150'
    'E42' = '  -0.72%  '
    'D43' = '0.895'
    'E43' = '  -2.89%  '
    'E44' = '  +4.41%  '
    'D45' = '5.68'
    'E45' = '  -0.32%  '
    'B46' = 'WhiteBITCoin'
    'C46' = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    'D46' = '23.52'
    'E46' = '  -1.05%  '
    'B47' = 'VeChain'
    'C47' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D47' = '0.0420'
    'E47' = '  +1.14%  '
    'B48' = 'Stacks'
    'C48' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D48' = '2.26'
    'E48' = '  +0.72%  '
    'D49' = '33.47'
    'E49' = '  -3.34%  '
    'E50' = '  -0.80%  '
    'E51' = '  +0.50%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
